$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3 : the "fichier d'entree" / "fichier resultat" scrolls + arrows
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)

# Existing shapes (before the edit): the "Parchemin vertical 79" scroll
# (id 80) and its "Fleche vers le bas 80" arrow (id 81), near the top of
# the slide, with no caption text on the scroll.
$scroll1 = $null
$arrow1 = $null
for ($i = 1; $i -le $s3.Shapes.Count; $i++) {
    $shp = $s3.Shapes.Item($i)
    if ($shp.Id -eq 80) { $scroll1 = $shp }
    if ($shp.Id -eq 81) { $arrow1 = $shp }
}

# Duplicate the pair so the new "Fichier Resultat" group keeps the exact
# same style (line/fill/effect/font refs) as the original.
$scroll2Range = $scroll1.Duplicate()
$scroll2 = $scroll2Range.Item(1)
$arrow2Range = $arrow1.Duplicate()
$arrow2 = $arrow2Range.Item(1)

# Resize / reposition the original scroll and give it its caption.
$scroll1.Left = 4788024 / 12700
$scroll1.Top = 908720 / 12700
$scroll1.Width = 1296144 / 12700
$scroll1.Height = 648072 / 12700
$scroll1.TextFrame.TextRange.Text = "Fichier d’entrée"
$scroll1.TextFrame.TextRange.Font.Size = 18

# The first arrow (id 81) keeps its original position/size - untouched.

# Position the duplicated scroll (new "Fichier Resultat") further down
# the slide and set its caption text.
$scroll2.Name = "Parchemin vertical 70"
$scroll2.Left = 4860032 / 12700
$scroll2.Top = 6021288 / 12700
$scroll2.Width = 1296144 / 12700
$scroll2.Height = 648072 / 12700
$scroll2.TextFrame.TextRange.Text = "Fichier Résultat"
$scroll2.TextFrame.TextRange.Font.Size = 18

# Position the duplicated arrow below the new scroll.
$arrow2.Name = "Fleche vers le bas 76"
$arrow2.Left = 5292080 / 12700
$arrow2.Top = 5445224 / 12700
$arrow2.Width = 216024 / 12700
$arrow2.Height = 504056 / 12700

# ---------------------------------------------------------------------
# Slide 4 : merge the "Etat " / "d'avancement" runs in the table header
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
for ($i = 1; $i -le $s4.Shapes.Count; $i++) {
    $shp = $s4.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $cell = $tbl.Cell($r, $c)
                $txt = $cell.Shape.TextFrame.TextRange.Text
                if ($txt -eq "Etat d’avancement") {
                    $cell.Shape.TextFrame.TextRange.Text = "Etat d’avancement"
                }
            }
        }
    }
}
